# Updated cryptos list on Tue Apr 11 23:56:33 UTC 2023 with GitHub Actions
# Refreshes the Price (D) and Volume(1h) (E) columns for each coin row.
# For D-column values that look like plain numbers, a leading apostrophe
# forces Excel to keep them as text (matching the original inlineStr cells
# instead of letting Excel coerce them to numeric), and the style is reset
# back to "Normal" afterwards so no stray number-format/quote-prefix style
# is left attached to the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '30.273.96'
$ws.Range('E2').Value = '  +2.05%  '
$ws.Range('D3').Value = '1.893.86'
$ws.Range('E3').Value = '  -0.89%  '
$ws.Range('E4').Value = '  -0.11%  '
$ws.Range('D5').Value = '''322.83'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +1.39%  '
$ws.Range('E6').Value = '  -0.09%  '
$ws.Range('D7').Value = '''0.5178'
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = '  -0.15%  '
$ws.Range('D8').Value = '''0.4017'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  +1.10%  '
$ws.Range('D9').Value = '''0.08431'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  -0.71%  '
$ws.Range('D10').Value = '''42.69'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  -0.47%  '
$ws.Range('E11').Value = '  -0.65%  '
$ws.Range('D12').Value = '''23.04'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  +10.34%  '
$ws.Range('D13').Value = '''6.431'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  +2.23%  '
$ws.Range('D14').Value = '1.892.40'
$ws.Range('E14').Value = '  -0.67%  '
$ws.Range('D15').Value = '''7.315'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  -0.41%  '
$ws.Range('D16').Value = '''1.002'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  -0.08%  '
$ws.Range('D17').Value = '''94.32'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  +0.43%  '
$ws.Range('E18').Value = '  -0.62%  '
$ws.Range('D19').Value = '''0.06653'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  -1.47%  '
$ws.Range('E20').Value = '  +1.40%  '
$ws.Range('D21').Value = '''1.000'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  -0.02%  '
$ws.Range('E22').Value = '  -1.14%  '
$ws.Range('D23').Value = '30.261.12'
$ws.Range('E23').Value = '  +1.98%  '
$ws.Range('E24').Value = '  +0.84%  '
$ws.Range('D25').Value = '''2.229'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  +0.94%  '
$ws.Range('D26').Value = '2.106.82'
$ws.Range('E26').Value = '  -0.84%  '
$ws.Range('D27').Value = '''21.51'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  +2.51%  '
$ws.Range('E28').Value = '  +1.75%  '
$ws.Range('D29').Value = '''2.334'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  -5.24%  '
$ws.Range('D30').Value = '''129.14'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  +0.47%  '
$ws.Range('D31').Value = '''1.086'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  +0.23%  '
$ws.Range('E32').Value = '  -0.64%  '
$ws.Range('D33').Value = '''6.118'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  -0.98%  '
$ws.Range('D34').Value = '''3.741'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  +1.74%  '
$ws.Range('D35').Value = '''0.02494'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  -0.21%  '
$ws.Range('D36').Value = '''0.06533'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  -1.34%  '
$ws.Range('D37').Value = '''5.343'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  +3.00%  '
$ws.Range('D38').Value = '''0.2202'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  -0.37%  '
$ws.Range('E39').Value = '  -1.85%  '
$ws.Range('D40').Value = '''8.797'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  -3.82%  '
$ws.Range('E41').Value = '  +3.09%  '
$ws.Range('E42').Value = '  -0.49%  '
$ws.Range('D43').Value = '''1.228'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  -1.18%  '
$ws.Range('D44').Value = '''0.6091'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  -0.77%  '
$ws.Range('D45').Value = '''13.29'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  +0.59%  '
$ws.Range('D46').Value = '''3.678'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  -0.52%  '
$ws.Range('D47').Value = '''2.056'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  -0.40%  '
$ws.Range('E48').Value = '  -0.19%  '
$ws.Range('D49').Value = '''124.75'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  +0.13%  '
$ws.Range('D50').Value = '''1.158'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  -2.58%  '
$ws.Range('D51').Value = '''79.20'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  +0.99%  '
